$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update nutrition values (columns H:M, rows 3-37) to be per-gram (divide by 100) ---
# and the scaling reference cell H48 (divide by 10, from 1000 -> 100)
$ws.Range("H3").Value = 0.5
$ws.Range("I3").Value = 0.00001
$ws.Range("J3").Value = 0.044000000000000004
$ws.Range("K3").Value = 0.22
$ws.Range("L3").Value = 0.21
$ws.Range("M3").Value = 0.0
$ws.Range("H4").Value = 0.002
$ws.Range("I4").Value = 0.00001
$ws.Range("J4").Value = 0.1
$ws.Range("K4").Value = 0.14
$ws.Range("L4").Value = 0.003
$ws.Range("M4").Value = 0.0
$ws.Range("H5").Value = 0.003
$ws.Range("I5").Value = 0.00001
$ws.Range("J5").Value = 0.12
$ws.Range("K5").Value = 0.23
$ws.Range("L5").Value = 0.011000000000000001
$ws.Range("M5").Value = 0.0
$ws.Range("H6").Value = 0.13
$ws.Range("I6").Value = 0.00068
$ws.Range("J6").Value = 0.0
$ws.Range("K6").Value = 0.0
$ws.Range("L6").Value = 0.19
$ws.Range("M6").Value = 0.00062
$ws.Range("H7").Value = 0.032
$ws.Range("I7").Value = 0.00508
$ws.Range("J7").Value = 0.059000000000000004
$ws.Range("K7").Value = 0.49
$ws.Range("L7").Value = 0.11
$ws.Range("M7").Value = 0.0
$ws.Range("H8").Value = 0.81
$ws.Range("I8").Value = 0.00643
$ws.Range("J8").Value = 0.001
$ws.Range("K8").Value = 0.001
$ws.Range("L8").Value = 0.009000000000000001
$ws.Range("M8").Value = 0.00215
$ws.Range("H9").Value = 0.001
$ws.Range("I9").Value = 0.00017999999999999998
$ws.Range("J9").Value = 0.032
$ws.Range("K9").Value = 0.057999999999999996
$ws.Range("L9").Value = 0.013000000000000001
$ws.Range("M9").Value = 0.0
$ws.Range("H10").Value = 0.002
$ws.Range("I10").Value = 0.0069
$ws.Range("J10").Value = 0.047
$ws.Range("K10").Value = 0.096
$ws.Range("L10").Value = 0.009000000000000001
$ws.Range("M10").Value = 0.0
$ws.Range("H11").Value = 0.31
$ws.Range("I11").Value = 0.00187
$ws.Range("J11").Value = 0.0
$ws.Range("K11").Value = 0.013999999999999999
$ws.Range("L11").Value = 0.27
$ws.Range("M11").Value = 0.0009299999999999999
$ws.Range("H12").Value = 0.18
$ws.Range("I12").Value = 0.00722
$ws.Range("J12").Value = 0.005
$ws.Range("K12").Value = 0.040999999999999995
$ws.Range("L12").Value = 0.18
$ws.Range("M12").Value = 0.0005600000000000001
$ws.Range("H13").Value = 0.027999999999999997
$ws.Range("I13").Value = 0.00246
$ws.Range("J13").Value = 0.04
$ws.Range("K13").Value = 0.23
$ws.Range("L13").Value = 0.071
$ws.Range("M13").Value = 0.0
$ws.Range("H14").Value = 0.047
$ws.Range("I14").Value = 0.00035000000000000005
$ws.Range("J14").Value = 0.006
$ws.Range("K14").Value = 0.74
$ws.Range("L14").Value = 0.094
$ws.Range("M14").Value = 0.0
$ws.Range("H15").Value = 0.002
$ws.Range("I15").Value = 0.00002
$ws.Range("J15").Value = 0.013999999999999999
$ws.Range("K15").Value = 0.022000000000000002
$ws.Range("L15").Value = 0.006
$ws.Range("M15").Value = 0.0
$ws.Range("H16").Value = 0.028999999999999998
$ws.Range("I16").Value = 0.00006
$ws.Range("J16").Value = 0.027000000000000003
$ws.Range("K16").Value = 0.73
$ws.Range("L16").Value = 0.14
$ws.Range("M16").Value = 0.0
$ws.Range("H17").Value = 0.15
$ws.Range("I17").Value = 0.00207
$ws.Range("J17").Value = 0.004
$ws.Range("K17").Value = 0.008
$ws.Range("L17").Value = 0.14
$ws.Range("M17").Value = 0.0040100000000000005
$ws.Range("H18").Value = 0.002
$ws.Range("I18").Value = 0.0004
$ws.Range("J18").Value = 0.43
$ws.Range("K18").Value = 0.64
$ws.Range("L18").Value = 0.006999999999999999
$ws.Range("M18").Value = 0.0
$ws.Range("H19").Value = 0.23
$ws.Range("I19").Value = 0.0005899999999999999
$ws.Range("J19").Value = 0.0
$ws.Range("K19").Value = 0.0
$ws.Range("L19").Value = 0.17
$ws.Range("M19").Value = 0.00073
$ws.Range("H20").Value = 0.011000000000000001
$ws.Range("I20").Value = 0.00006
$ws.Range("J20").Value = 0.02
$ws.Range("K20").Value = 0.63
$ws.Range("L20").Value = 0.25
$ws.Range("M20").Value = 0.0
$ws.Range("H21").Value = 0.002
$ws.Range("I21").Value = 0.00028000000000000003
$ws.Range("J21").Value = 0.008
$ws.Range("K21").Value = 0.028999999999999998
$ws.Range("L21").Value = 0.013999999999999999
$ws.Range("M21").Value = 0.0
$ws.Range("H22").Value = 0.035
$ws.Range("I22").Value = 0.00003
$ws.Range("J22").Value = 0.045
$ws.Range("K22").Value = 0.045
$ws.Range("L22").Value = 0.031
$ws.Range("M22").Value = 0.00014000000000000001
$ws.Range("H23").Value = 0.15
$ws.Range("I23").Value = 0.015560000000000001
$ws.Range("J23").Value = 0.005
$ws.Range("K23").Value = 0.038
$ws.Range("L23").Value = 0.01
$ws.Range("M23").Value = 0.0
$ws.Range("H24").Value = 1.0
$ws.Range("I24").Value = 0.00002
$ws.Range("J24").Value = 0.0
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.0
$ws.Range("M24").Value = 0.0
$ws.Range("H25").Value = 0.001
$ws.Range("I25").Value = 0.00004
$ws.Range("J25").Value = 0.042
$ws.Range("K25").Value = 0.09300000000000001
$ws.Range("L25").Value = 0.011000000000000001
$ws.Range("M25").Value = 0.0
$ws.Range("H26").Value = 0.002
$ws.Range("I26").Value = 0.0
$ws.Range("J26").Value = 0.091
$ws.Range("K26").Value = 0.12
$ws.Range("L26").Value = 0.006999999999999999
$ws.Range("M26").Value = 0.0
$ws.Range("H27").Value = 1.0
$ws.Range("I27").Value = 0.0
$ws.Range("J27").Value = 0.0
$ws.Range("K27").Value = 0.0
$ws.Range("L27").Value = 0.0
$ws.Range("M27").Value = 0.0
$ws.Range("H28").Value = 0.004
$ws.Range("I28").Value = 0.00005
$ws.Range("J28").Value = 0.057
$ws.Range("K28").Value = 0.14
$ws.Range("L28").Value = 0.054000000000000006
$ws.Range("M28").Value = 0.0
$ws.Range("H29").Value = 0.21
$ws.Range("I29").Value = 0.0005600000000000001
$ws.Range("J29").Value = 0.0
$ws.Range("K29").Value = 0.0
$ws.Range("L29").Value = 0.17
$ws.Range("M29").Value = 0.0007199999999999999
$ws.Range("H30").Value = 0.001
$ws.Range("I30").Value = 0.00006
$ws.Range("J30").Value = 0.008
$ws.Range("K30").Value = 0.17
$ws.Range("L30").Value = 0.021
$ws.Range("M30").Value = 0.0
$ws.Range("H31").Value = 0.023
$ws.Range("I31").Value = 0.0005
$ws.Range("J31").Value = 0.0
$ws.Range("K31").Value = 0.0
$ws.Range("L31").Value = 0.22
$ws.Range("M31").Value = 0.0008100000000000001
$ws.Range("H32").Value = 0.011000000000000001
$ws.Range("I32").Value = 0.00231
$ws.Range("J32").Value = 0.038
$ws.Range("K32").Value = 0.21
$ws.Range("L32").Value = 0.08
$ws.Range("M32").Value = 0.0
$ws.Range("H33").Value = 0.032
$ws.Range("I33").Value = 0.00005
$ws.Range("J33").Value = 0.006999999999999999
$ws.Range("K33").Value = 0.76
$ws.Range("L33").Value = 0.075
$ws.Range("M33").Value = 0.0
$ws.Range("H34").Value = 0.01
$ws.Range("I34").Value = 0.00039
$ws.Range("J34").Value = 0.053
$ws.Range("K34").Value = 0.092
$ws.Range("L34").Value = 0.003
$ws.Range("M34").Value = 0.0
$ws.Range("H35").Value = 0.2
$ws.Range("I35").Value = 0.00016
$ws.Range("J35").Value = 0.027000000000000003
$ws.Range("K35").Value = 0.08900000000000001
$ws.Range("L35").Value = 0.19
$ws.Range("M35").Value = 0.0
$ws.Range("H36").Value = 0.003
$ws.Range("I36").Value = 0.0018599999999999999
$ws.Range("J36").Value = 0.044000000000000004
$ws.Range("K36").Value = 0.073
$ws.Range("L36").Value = 0.016
$ws.Range("M36").Value = 0.0
$ws.Range("H37").Value = 0.019
$ws.Range("I37").Value = 0.0011799999999999998
$ws.Range("J37").Value = 0.001
$ws.Range("K37").Value = 0.001
$ws.Range("L37").Value = 0.23
$ws.Range("M37").Value = 0.00067
$ws.Range("H48").Value = 100.0

# --- 2. Expand the AutoFilter range from A1:G31 to A1:M37 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:M37").AutoFilter()

# --- 3. Fix up the hidden _xlnm._FilterDatabase defined names so they track the
#        same history pattern Excel/Calc produces when the AutoFilter range changes:
#        the previous entries shift their ranges and a new one is appended. ---
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$M`$37"
$wb.Names.Item(2).RefersTo = "=Sheet1!`$A`$1:`$G`$31"
$wb.Names.Item(3).RefersTo = "=Sheet1!`$A`$1:`$M`$37"
$wb.Names.Item(4).RefersTo = "=Sheet1!`$A`$1:`$G`$31"
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0", "=Sheet1!`$A`$1:`$D`$23")

# --- 4. Move the active selection from D3 to O26 ---
$ws.Range("O26").Select()

Write-Host "edit complete"
